# Update the Tasks sheet: the list of generic "Task N" placeholder rows is
# replaced with the actual imaging-pipeline step names, and the two
# trailing "Kidney Volumes" / "AAA" rows are cleared out (no longer used).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("A2").Value = "Kidney Area"
$ws.Range("A5").Value = "Liver Area"
$ws.Range("A3").Value = "Kidney Vol."
$ws.Range("A4").Value = "Kidney Conc."
$ws.Range("A6").Value = "Liver Vol."
$ws.Range("A7").Value = "Liver Conc."
$ws.Range("A8").Value = "Decay Rate"

# Rows 9 and 10 no longer hold data - clear their contents, keeping the
# existing (unstriped) row formatting in place.
$ws.Range("A9:B10").ClearContents()

# Move the active selection to A8, matching where the editor left off.
$ws.Range("A8").Select()
